# Jasmine.pptx edit: fix slide 5 title run-split and append a new
# "Spies" slide (Title and Content layout) at the end of the deck.

$p = $ppt.ActivePresentation

# --- 1. Slide 5: merge the two title runs into a single run ---------------
# The title currently holds two runs ("Our " / "playground ... car") whose
# concatenation already equals the desired text, so re-assigning the same
# string directly only rewrites the changed tail (engine keeps unchanged
# run prefixes intact). Bounce through an unrelated value first so the
# whole paragraph gets rebuilt as one fresh run, then set the real text.
$s5 = $p.Slides.Item(5)
$s5Title = $s5.Shapes.Item(1).TextFrame.TextRange
$s5Title.Text = " "
$s5Title2 = $s5.Shapes.Item(1).TextFrame.TextRange
$s5Title2.Text = "Our playground " + [char]0x2013 + " a car"

# --- 2. Add the new "Spies" slide at the end (position 6) ------------------
$newSlide = $p.Slides.Add(6, 2)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Spies"

$apostrophe = [char]0x2019
$bodyText = "Useful when:`r" + `
    "Testing interactions`r" + `
    "Internal behavior (not a good practice)`r" + `
    "Or a functionality needed doesn" + $apostrophe + "t exist yet.`r"

$contentTr = $newSlide.Shapes.Item(2).TextFrame.TextRange
$contentTr.Text = $bodyText

# Paragraphs 2-4 ("Testing interactions" .. "...doesn't exist yet.") sit one
# indent level down (lvl=1 / IndentLevel=2).
$paraIndex = 0
ForEach ($para in $contentTr.Paragraphs()) {
    $paraIndex = $paraIndex + 1
    if ($paraIndex -gt 1) {
        $para.IndentLevel = 2
    }
}

# The trailing `\r` leaves an implicit 5th (empty) paragraph that isn't
# enumerated by Paragraphs() yet, but is addressable/indentable directly.
$trailingPara = $contentTr.Paragraphs(5, 1)
ForEach ($para in $trailingPara) {
    $para.IndentLevel = 2
}
